# Applies the "Gilds added + Clean up + Unittests" edit to the Checklist workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update individual score cells (column B) with newly-graded values.
$ws.Range("B14").Value = 1
$ws.Range("B36").Value = 3
$ws.Range("B41").Value = 3
$ws.Range("B42").Value = 0
$ws.Range("B48").Value = 2
$ws.Range("B49").Value = 2
$ws.Range("B50").Value = 2
$ws.Range("B51").Value = 1
$ws.Range("B52").Value = 1

# Recalculate so the B54 total formula reflects the new scores.
$excel.Calculate()

# Update the view state: scroll so row 22 is at the top, and select I40.
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("I40").Select()
